$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row containing "io" / "E" (row 5) — rows below shift up.
$ws.Rows.Item(5).Delete()

# Restore the selection Excel leaves behind after this edit.
$ws.Range("F15").Select()
